$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "g"

$ws.Range("B7").Select()
